$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their text representation (avoid Excel auto-numeric conversion)
$textCells = @(
    "D2", "D3", "D5", "D9", "D10", "D12", "D13", "D14", "D16", "D17", "D19", "D20", "D23", "D24", "D25", "D29", "D32", "D33", "D37", "D38", "D39", "D40", "D45", "D48", "D49", "D50", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values per the diff
$ws.Range('D2').Value = '26.142.62'
$ws.Range('E2').Value = '  -0.41%  '
$ws.Range('D3').Value = '1.583.84'
$ws.Range('E3').Value = '  -0.04%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '211.24'
$ws.Range('E5').Value = '  +1.04%  '
$ws.Range('E6').Value = '  +0.22%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = '0.0606'
$ws.Range('E9').Value = '  -0.88%  '
$ws.Range('D10').Value = '19.19'
$ws.Range('E10').Value = '  -2.02%  '
$ws.Range('E11').Value = '  +0.23%  '
$ws.Range('D12').Value = '1.807.00'
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '4.01'
$ws.Range('E13').Value = '  -1.43%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.565.60'
$ws.Range('E14').Value = '  -1.71%  '
$ws.Range('E15').Value = '  -0.09%  '
$ws.Range('D16').Value = '64.00'
$ws.Range('E16').Value = '  -1.04%  '
$ws.Range('D17').Value = '26.141.57'
$ws.Range('E17').Value = '  -0.42%  '
$ws.Range('E18').Value = '  -0.64%  '
$ws.Range('D19').Value = '7.33'
$ws.Range('E19').Value = '  -0.76%  '
$ws.Range('D20').Value = '213.14'
$ws.Range('E20').Value = '  +0.30%  '
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('E22').Value = '  -0.69%  '
$ws.Range('D23').Value = '2.17'
$ws.Range('E23').Value = '  -0.39%  '
$ws.Range('D24').Value = '8.93'
$ws.Range('E24').Value = '  +0.75%  '
$ws.Range('D25').Value = '143.83'
$ws.Range('E25').Value = '  -0.55%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('E27').Value = '  -0.83%  '
$ws.Range('E28').Value = '  -0.73%  '
$ws.Range('D29').Value = '15.12'
$ws.Range('E29').Value = '  -1.17%  '
$ws.Range('E30').Value = '  -2.01%  '
$ws.Range('E31').Value = '  +0.56%  '
$ws.Range('D32').Value = '3.19'
$ws.Range('E32').Value = '  -1.43%  '
$ws.Range('D33').Value = '1.338.51'
$ws.Range('E33').Value = '  +3.94%  '
$ws.Range('E34').Value = '  -2.09%  '
$ws.Range('E35').Value = '  -0.11%  '
$ws.Range('E36').Value = '  -1.30%  '
$ws.Range('D37').Value = '0.580'
$ws.Range('E37').Value = '  -3.98%  '
$ws.Range('D38').Value = '0.0167'
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('D39').Value = '0.816'
$ws.Range('E39').Value = '  +0.41%  '
$ws.Range('D40').Value = '5.79'
$ws.Range('E40').Value = '  +2.71%  '
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('E42').Value = '  -16.84%  '
$ws.Range('E43').Value = '  +0.33%  '
$ws.Range('E44').Value = '  -0.44%  '
$ws.Range('D45').Value = '1.719.36'
$ws.Range('E46').Value = '  -2.73%  '
$ws.Range('E47').Value = '  -3.09%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0103'
$ws.Range('E48').Value = '  +5.72%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '1.47'
$ws.Range('E49').Value = '  -1.87%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '0.0982'
$ws.Range('E50').Value = '  -1.62%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '0.0500'
$ws.Range('E51').Value = '  -1.22%  '
